$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (legmaxROM updates)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) - meanEMG updated; C2 and E2 cleared
$ws.Range("B2").Value = 17.69951367420683
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 26.021814370035074
$ws.Range("E2").ClearContents()

# Row 3 (STR) - meanEMG updated
$ws.Range("B3").Value = 14.471734435433772
$ws.Range("C3").Value = -10.570762868888698
$ws.Range("D3").Value = 26.152174988794961
$ws.Range("E3").Value = -4.6204314984218495

# Update selection to match new sqref B1:E3
$ws.Range("B1:E3").Select()
